# Apply the 2017-01-31 EIA Table A.6.B monthly update (Year-to-Date through
# November 2016, replacing the prior October 2016 edition).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the subtitle text (row 2) to reflect the new reporting month.
$ws.Range("A2").Value = "by End-Use Sector, Census Division, and State, Year-to-Date through November 2016"

# Update the relative-standard-error data values that changed between the
# October 2016 and November 2016 runs.
$ws.Range("B4").Value  = 0.15
$ws.Range("B5").Value  = 0.14000000000000001
$ws.Range("B6").Value  = 0.21
$ws.Range("F6").Value  = 0.44
$ws.Range("B7").Value  = 0.33
$ws.Range("B8").Value  = 0.22
$ws.Range("C12").Value = 0.22
$ws.Range("B14").Value = 0.12
$ws.Range("C14").Value = 0.28999999999999998
$ws.Range("F14").Value = 0.17
$ws.Range("B15").Value = 0.14000000000000001
$ws.Range("B16").Value = 0.22
$ws.Range("B17").Value = 0.39
$ws.Range("B18").Value = 0.2
$ws.Range("B19").Value = 0.24
$ws.Range("B20").Value = 0.39
$ws.Range("B21").Value = 0.25
$ws.Range("B25").Value = 0.41
$ws.Range("B29").Value = 0.19
$ws.Range("C29").Value = 0.19
$ws.Range("B30").Value = 0.42
$ws.Range("B32").Value = 0.24
$ws.Range("C32").Value = 0.32
$ws.Range("B34").Value = 0.14000000000000001
$ws.Range("B35").Value = 0.41
$ws.Range("C35").Value = 0.45
$ws.Range("F35").Value = 0.43
$ws.Range("B37").Value = 0.32
$ws.Range("B38").Value = 0.1
$ws.Range("C38").Value = 0.38
$ws.Range("B41").Value = 0.48
$ws.Range("B43").Value = 0.34
$ws.Range("B44").Value = 0.28000000000000003
$ws.Range("C44").Value = 0.26
$ws.Range("B48").Value = 0.3
$ws.Range("B52").Value = 0.45
$ws.Range("B54").Value = 0.19
$ws.Range("B58").Value = 0.14000000000000001
$ws.Range("B60").Value = 0.46
$ws.Range("D60").Value = 6
$ws.Range("B61").Value = 0.37
$ws.Range("B65").Value = 0.1

$wb.Save()
